$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.906979
$ws.Range("H2").Value = 107.720937
$ws.Range("I2").Value = 0.6107087147789413
$ws.Range("J2").Value = 0.6107087147789412
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 135.00078553553
$ws.Range("R2").Value = 1215.00706981977
$ws.Range("S2").Value = 0.04175731235202055
$ws.Range("T2").Value = 0.04175731235202054

# Row 3
$ws.Range("G3").Value = 35.906979
$ws.Range("H3").Value = 107.720937
$ws.Range("I3").Value = 0.6107087147789413
$ws.Range("J3").Value = 0.6107087147789412
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("Q3").Value = 1286.167796088809
$ws.Range("R3").Value = 11575.51016479928
$ws.Range("S3").Value = 0.3978266510475637
$ws.Range("T3").Value = 0.3978266510475635

# Row 4
$ws.Range("G4").Value = 35.906979
$ws.Range("H4").Value = 107.720937
$ws.Range("I4").Value = 0.6107087147789413
$ws.Range("J4").Value = 0.6107087147789412
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("Q4").Value = 553.24384064837
$ws.Range("R4").Value = 4979.19456583533
$ws.Range("S4").Value = 0.171124751379357
$ws.Range("T4").Value = 0.171124751379357

# Row 5
$ws.Range("I5").Value = 0.2899643113254147
$ws.Range("J5").Value = 0.2899643113254147
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 64.09833175603111
$ws.Range("R5").Value = 576.8849858042801
$ws.Range("S5").Value = 0.01982635915607764
$ws.Range("T5").Value = 0.01982635915607764

# Row 6
$ws.Range("I6").Value = 0.2899643113254147
$ws.Range("J6").Value = 0.2899643113254147
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("S6").Value = 0.1888879724594371
$ws.Range("T6").Value = 0.188887972459437

# Row 7
$ws.Range("I7").Value = 0.2899643113254147
$ws.Range("J7").Value = 0.2899643113254147
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("S7").Value = 0.08124997970989999
$ws.Range("T7").Value = 0.08124997970989996

# Row 8
$ws.Range("I8").Value = 0.09932697389564409
$ws.Range("J8").Value = 0.09932697389564407
$ws.Range("M8").Value = 3.759736666666667
$ws.Range("N8").Value = 11.27921
$ws.Range("O8").Value = 0.0683751702595819
$ws.Range("P8").Value = 0.06837517025958188
$ws.Range("Q8").Value = 21.95681701649334
$ws.Range("R8").Value = 197.61135314844
$ws.Range("S8").Value = 0.006791498751483711
$ws.Range("T8").Value = 0.006791498751483709

# Row 9
$ws.Range("I9").Value = 0.09932697389564409
$ws.Range("J9").Value = 0.09932697389564407
$ws.Range("O9").Value = 0.6514180024294648
$ws.Range("P9").Value = 0.6514180024294647
$ws.Range("S9").Value = 0.06470337892246407
$ws.Range("T9").Value = 0.06470337892246404

# Row 10
$ws.Range("I10").Value = 0.09932697389564409
$ws.Range("J10").Value = 0.09932697389564407
$ws.Range("O10").Value = 0.2802068273109533
$ws.Range("P10").Value = 0.2802068273109533
$ws.Range("S10").Value = 0.02783209622169631
$ws.Range("T10").Value = 0.0278320962216963
